# Update the "Förändrad" (Changed) date column (C) for rows 2 through 20
# from 2023-09-20 (serial 45189) to 2023-09-21 (serial 45190).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
